$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 991864.0699999999
$ws.Range("E2").Value = 645673.29
$ws.Range("H2").Value = 28880
$ws.Range("I2").Value = 28880
$ws.Range("J2").Value = 28880
$ws.Range("K2").Value = 981451.0699999999
$ws.Range("D3").Value = 812097.11
$ws.Range("E3").Value = 651365.6800000001
$ws.Range("G3").Value = 23651.42
$ws.Range("H3").Value = 23641
$ws.Range("I3").Value = 23652.02
$ws.Range("J3").Value = 23556.64
$ws.Range("K3").Value = 806962.11
$ws.Range("L3").Value = -174488.96
$ws.Range("M3").Value = -33.98
$ws.Range("N3").Value = 5692.39
$ws.Range("Q3").Value = -174488.96
$ws.Range("D4").Value = 893206.35
$ws.Range("E4").Value = 662412.49
$ws.Range("G4").Value = 26007.39
$ws.Range("H4").Value = 26007.39
$ws.Range("I4").Value = 26007.39
$ws.Range("J4").Value = 26007.39
$ws.Range("K4").Value = 885939.35
$ws.Range("L4").Value = -95511.72
$ws.Range("M4").Value = -13.14
$ws.Range("N4").Value = 16739.2
$ws.Range("Q4").Value = -95511.72
$ws.Range("D5").Value = 928237.5600000001
$ws.Range("E5").Value = 662412.49
$ws.Range("G5").Value = 27027.39
$ws.Range("H5").Value = 27027.39
$ws.Range("I5").Value = 27027.39
$ws.Range("J5").Value = 27027.39
$ws.Range("K5").Value = 920489.5600000001
$ws.Range("L5").Value = -60961.51
$ws.Range("M5").Value = -7.87
$ws.Range("N5").Value = 16739.2
$ws.Range("Q5").Value = -60961.51
$ws.Range("D6").Value = 938540.86
$ws.Range("E6").Value = 662412.49
$ws.Range("G6").Value = 27327.39
$ws.Range("H6").Value = 27327.39
$ws.Range("I6").Value = 27327.39
$ws.Range("J6").Value = 27327.39
$ws.Range("K6").Value = 930259.86
$ws.Range("L6").Value = -51191.22
$ws.Range("M6").Value = -6.18
$ws.Range("N6").Value = 16739.2
$ws.Range("Q6").Value = -51191.22
$ws.Range("D7").Value = 961191.76
$ws.Range("E7").Value = 660335.48
$ws.Range("G7").Value = 27986.92
$ws.Range("H7").Value = 27986.92
$ws.Range("I7").Value = 27986.92
$ws.Range("J7").Value = 27986.92
$ws.Range("K7").Value = 952377.76
$ws.Range("L7").Value = -29073.31
$ws.Range("M7").Value = -3.3
$ws.Range("N7").Value = 14662.2
$ws.Range("Q7").Value = -29073.31
$ws.Range("D8").Value = 967492.05
$ws.Range("E8").Value = 656332.48
$ws.Range("G8").Value = 28170.36
$ws.Range("H8").Value = 28170.36
$ws.Range("I8").Value = 28170.36
$ws.Range("J8").Value = 28170.36
$ws.Range("K8").Value = 958145.05
$ws.Range("L8").Value = -23306.02
$ws.Range("M8").Value = -2.49
$ws.Range("N8").Value = 10659.19
$ws.Range("Q8").Value = -23306.02
$ws.Range("D9").Value = 988003.22
$ws.Range("E9").Value = 652115.73
$ws.Range("G9").Value = 28767.58
$ws.Range("H9").Value = 28767.58
$ws.Range("I9").Value = 28767.58
$ws.Range("J9").Value = 28767.58
$ws.Range("K9").Value = 978123.22
$ws.Range("L9").Value = -3327.85
$ws.Range("M9").Value = -0.34
$ws.Range("N9").Value = 6442.44
$ws.Range("Q9").Value = -3327.85
$ws.Range("D10").Value = 994358.84
$ws.Range("E10").Value = 637864.76
$ws.Range("G10").Value = 28952.64
$ws.Range("H10").Value = 28952.64
$ws.Range("I10").Value = 28952.64
$ws.Range("J10").Value = 28952.64
$ws.Range("K10").Value = 983375.84
$ws.Range("L10").Value = 1924.77
$ws.Range("M10").Value = 0.18
$ws.Range("N10").Value = -7808.53
$ws.Range("Q10").Value = 1924.77
$ws.Range("D11").Value = 1011810.25
$ws.Range("E11").Value = 630588.26
$ws.Range("G11").Value = 29460.77
$ws.Range("H11").Value = 29460.77
$ws.Range("I11").Value = 29460.77
$ws.Range("J11").Value = 29460.77
$ws.Range("K11").Value = 1000253.25
$ws.Range("L11").Value = 18802.18
$ws.Range("M11").Value = 1.63
$ws.Range("N11").Value = -15085.03
$ws.Range("Q11").Value = 18802.18
$ws.Range("D12").Value = 1011867.9
$ws.Range("E12").Value = 620342.61
$ws.Range("G12").Value = 29462.45
$ws.Range("H12").Value = 29462.45
$ws.Range("I12").Value = 29462.45
$ws.Range("J12").Value = 29462.45
$ws.Range("K12").Value = 999740.9
$ws.Range("L12").Value = 18289.83
$ws.Range("M12").Value = 1.51
$ws.Range("N12").Value = -25330.68
$ws.Range("Q12").Value = 18289.83
$ws.Range("D13").Value = 1022544.64
$ws.Range("E13").Value = 606291.4300000001
$ws.Range("G13").Value = 29773.32
$ws.Range("H13").Value = 29773.32
$ws.Range("I13").Value = 29773.32
$ws.Range("J13").Value = 29773.32
$ws.Range("K13").Value = 1009843.64
$ws.Range("L13").Value = 28392.57
$ws.Range("M13").Value = 2.24
$ws.Range("N13").Value = -39381.85
$ws.Range("Q13").Value = 28392.57
$ws.Range("D14").Value = 1021649.75
$ws.Range("E14").Value = 595093.24
$ws.Range("G14").Value = 29747.27
$ws.Range("H14").Value = 29747.27
$ws.Range("I14").Value = 29747.27
$ws.Range("J14").Value = 29747.27
$ws.Range("K14").Value = 1008378.75
$ws.Range("L14").Value = 26927.67
$ws.Range("M14").Value = 2.03
$ws.Range("N14").Value = -50580.04
$ws.Range("Q14").Value = 26927.67
$ws.Range("D15").Value = 1040069.53
$ws.Range("E15").Value = 578481.8199999999
$ws.Range("G15").Value = 30283.59
$ws.Range("H15").Value = 30283.59
$ws.Range("I15").Value = 30283.59
$ws.Range("J15").Value = 30283.59
$ws.Range("K15").Value = 1026278.53
$ws.Range("L15").Value = 44827.46
$ws.Range("M15").Value = 3.25
$ws.Range("N15").Value = -67191.47
$ws.Range("Q15").Value = 44827.46
$ws.Range("D16").Value = 1046181.45
$ws.Range("E16").Value = 514531.32
$ws.Range("G16").Value = 30461.55
$ws.Range("H16").Value = 30461.55
$ws.Range("I16").Value = 30461.55
$ws.Range("J16").Value = 30461.55
$ws.Range("K16").Value = 1030102.45
$ws.Range("L16").Value = 48651.38
$ws.Range("M16").Value = 3.03
$ws.Range("N16").Value = -131141.97
$ws.Range("Q16").Value = 48651.38
$ws.Range("D17").Value = 966600.16
$ws.Range("E17").Value = 619023.71
$ws.Range("H17").Value = 27521.81
$ws.Range("I17").Value = 27993.65
$ws.Range("J17").Value = 26724.67
$ws.Range("K17").Value = 956187.16
$ws.Range("L17").Value = -25263.92
$ws.Range("M17").Value = -2.43
$ws.Range("N17").Value = -26649.58
$ws.Range("Q17").Value = -26649.58
